# Control.xlsx update: add two new certification rows (45, 46) with their
# hyperlinks, matching the "Atualizacao e criacao de arquivos." commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add rows 45 and 46, cloning the formatting of row 44 (the last
#     existing data row) so the new cells pick up the same styles
#     (s="2" text, s="3" date, s="4" hyperlink) instead of the workbook
#     default style. ---
$ws.Range("B44:I44").Copy()
$ws.Range("B45:I45").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B46:I46").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 45: AWS Certified Cloud Practitioner --------------------------
# Fill the plain / reused-string values first.
$ws.Range("B45").Value = "Aws"
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 45187
$ws.Range("G45").Value = "Ok"
$ws.Range("H45").Value = "Ok"
$ws.Range("I45").Value = 45239

# --- Row 46: Governanca de Dados (Data Science Academy) ---------------
$ws.Range("B46").Value = "Data Science Academy"
$ws.Range("C46").Value = "Governança de Dados"
$ws.Range("D46").Value = 8
$ws.Range("E46").Value = 45230
$ws.Range("G46").Value = "Ok"
$ws.Range("H46").Value = "Ok"
$ws.Range("I46").Value = 45239
$ws.Hyperlinks.Add($ws.Range("F46"), "https://mycourse.app/MVxqobEMGb9Dp9LW6") | Out-Null
# Hyperlinks.Add applies the built-in "Hyperlink" style (11pt) instead of
# keeping the 9pt hyperlink style already on the row; re-apply the row's
# format (value is untouched, PasteSpecial formats-only) to restore it.
$ws.Range("F44").Copy()
$ws.Range("F46").PasteSpecial(-4122)   # xlPasteFormats

# Finish row 45's remaining new strings (kept after row 46 so the shared
# string table grows in the same order as the authored workbook).
$ws.Range("C45").Value = "AWS Certified Cloud Practitioner"
$ws.Hyperlinks.Add($ws.Range("F45"), "https://www.credly.com/badges/38cdec31-4934-47a5-8da4-a51b01640397/linked_in_profile") | Out-Null
$ws.Range("F44").Copy()
$ws.Range("F45").PasteSpecial(-4122)   # xlPasteFormats

# --- Update the view state: selection moved to B48, no more scrolled
#     topLeftCell override. ---
$ws.Range("B48").Select() | Out-Null
